# Auto-generated Excel COM-interop script that applies the crypto price
# table refresh described in the commit diff for cryptos.xlsx.
#
# For each changed row it updates:
#   - B/C (coin name / coinranking.com link) where the row content changed
#   - D   (Price) - forced to Text storage so Excel does not reinterpret
#         values such as "0.664" or "43.745.98" as numbers/dates
#   - E   (Volume(1h) percentage string, incl. surrounding spaces)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain text values (e.g. "0.664", "43.745.98")
# that Excel would otherwise auto-convert to numbers when assigned via
# .Value. Temporarily mark the column as Text, write the values, then
# restore the default "Normal" style so the cells keep no explicit
# number-format override (matching the original workbook).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.745.98"
$ws.Range("E2").Value = "  -0.47%  "

# Row 3
$ws.Range("D3").Value = "2.337.50"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "238.85"
$ws.Range("E5").Value = "  -1.15%  "

# Row 6
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -3.60%  "

# Row 7
$ws.Range("D7").Value = "71.75"
$ws.Range("E7").Value = "  -6.46%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -7.57%  "

# Row 10
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  -3.80%  "

# Row 11
$ws.Range("D11").Value = "57.64"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12
$ws.Range("D12").Value = "32.23"
$ws.Range("E12").Value = "  -2.37%  "

# Row 13
$ws.Range("E13").Value = "  -0.72%  "

# Row 14
$ws.Range("D14").Value = "7.09"
$ws.Range("E14").Value = "  -6.61%  "

# Row 15
$ws.Range("D15").Value = "2.685.64"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").Value = "16.12"
$ws.Range("E16").Value = "  -4.70%  "

# Row 17
$ws.Range("E17").Value = "  -3.48%  "

# Row 18
$ws.Range("D18").Value = "2.343.49"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
$ws.Range("D19").Value = "43.602.58"
$ws.Range("E19").Value = "  -0.78%  "

# Row 20
$ws.Range("E20").Value = "  -2.09%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "78.16"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").Value = "251.90"
$ws.Range("E23").Value = "  -2.32%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("E25").Value = "  +7.98%  "

# Row 26
$ws.Range("E26").Value = "  +1.54%  "

# Row 27
$ws.Range("D27").Value = "2.48"
$ws.Range("E27").Value = "  -2.36%  "

# Row 28
$ws.Range("D28").Value = "10.31"
$ws.Range("E28").Value = "  -8.11%  "

# Row 29
$ws.Range("E29").Value = "  -0.95%  "

# Row 30
$ws.Range("D30").Value = "175.03"
$ws.Range("E30").Value = "  -0.49%  "

# Row 31
$ws.Range("D31").Value = "22.10"
$ws.Range("E31").Value = "  -4.49%  "

# Row 32
$ws.Range("E32").Value = "  -2.55%  "

# Row 33
$ws.Range("E33").Value = "  -0.41%  "

# Row 34
$ws.Range("D34").Value = "0.0732"
$ws.Range("E34").Value = "  -2.44%  "

# Row 35
$ws.Range("D35").Value = "5.06"
$ws.Range("E35").Value = "  -5.00%  "

# Row 36
$ws.Range("E36").Value = "  -0.56%  "

# Row 37
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  -2.45%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -3.92%  "

# Row 39
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "6.33"
$ws.Range("E39").Value = "  -2.68%  "

# Row 40
$ws.Range("D40").Value = "0.0268"
$ws.Range("E40").Value = "  -2.72%  "

# Row 41
$ws.Range("D41").Value = "5.26"
$ws.Range("E41").Value = "  +16.79%  "

# Row 42
$ws.Range("D42").Value = "64.44"
$ws.Range("E42").Value = "  +17.48%  "

# Row 43
$ws.Range("E43").Value = "  +1.90%  "

# Row 44
$ws.Range("D44").Value = "0.106"
$ws.Range("E44").Value = "  +5.51%  "

# Row 45
$ws.Range("D45").Value = "18.81"
$ws.Range("E45").Value = "  -1.12%  "

# Row 46
$ws.Range("E46").Value = "  -4.13%  "

# Row 47
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("E48").Value = "  -3.76%  "

# Row 49
$ws.Range("E49").Value = "  -3.41%  "

# Row 50
$ws.Range("E50").Value = "  -5.17%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "97.28"
$ws.Range("E51").Value = "  -5.01%  "

# Restore the Price column cells to the default style (removes the
# temporary Text number-format override applied above).
$priceRange.Style = "Normal"
